$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell A2: "Emre Abale" -> "Rob Oudman"
$ws.Range("A2").Value = "Rob Oudman"

# Update the active selection to A2 (matches sheetView selection change)
$ws.Range("A2").Select()
